# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1) Metadata sheet: bump the StructureDefinition "Date" value.
# 2) Elements sheet: the two mapping columns (AK = "Mapping: RIM Mapping",
#    AL = "Mapping: Spécification métier vers l'extension ROR
#    OrganizationCreationDate") are swapped - the French business-mapping
#    column now comes first, RIM Mapping second. Swap header text, the
#    per-row values, and the column widths that went with them.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata!B8 : Date -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2) Elements : swap columns AK (37) and AL (38) ------------------------
$els = $wb.Worksheets.Item("Elements")

# Remember the current contents of both columns, row by row, before
# overwriting anything. (`.Text` is used for the read-back because this
# host's `Range.Value` getter isn't wired up; `.Text` is.)
$akHeader = $els.Range("AK1").Text
$alHeader = $els.Range("AL1").Text
$ak3 = $els.Range("AK3").Text
$al3 = $els.Range("AL3").Text
$ak4 = $els.Range("AK4").Text
$al4 = $els.Range("AL4").Text
$ak5 = $els.Range("AK5").Text
$al5 = $els.Range("AL5").Text
$ak6 = $els.Range("AK6").Text
$al6 = $els.Range("AL6").Text

# Write back with AK/AL exchanged.
$els.Range("AK1").Value = $alHeader
$els.Range("AL1").Value = $akHeader

$els.Range("AK3").Value = $al3
$els.Range("AL3").Value = $ak3

$els.Range("AK4").Value = $al4
$els.Range("AL4").Value = $ak4

$els.Range("AK5").Value = $al5
$els.Range("AL5").Value = $ak5

$els.Range("AK6").Value = $al6
$els.Range("AL6").Value = $ak6

# Column widths followed the data: AK becomes the wide "Spécification
# métier" column, AL becomes the narrower "RIM Mapping" column.
$els.Columns.Item(37).ColumnWidth = 81.16666666666667
$els.Columns.Item(38).ColumnWidth = 24.166666666666664
